# Add the new "2022-Q3" quarterly sheet right after "总计", shifting the
# existing quarter sheets (2022-Q2 .. 2021-Q1) one position to the right,
# and insert the corresponding summary row into "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q3" right after "总计" (position 2).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item(1))
$newSheet.Name = "2022-Q3"

# Copy the header-row formatting (bold / bordered style) and the data-row
# formatting (column A index style) from the existing "2022-Q2" sheet so the
# new sheet reuses the same style indices instead of minting new ones.
$srcHeader = $wb.Worksheets.Item("2022-Q2")
$srcHeader.Range("A1:H1").Copy()
$dstHeader = $wb.Worksheets.Item("2022-Q3")
$dstHeader.Range("A1:H1").PasteSpecial(-4122)   # xlPasteFormats

$srcData = $wb.Worksheets.Item("2022-Q2")
$srcData.Range("A2:H2").Copy()
$dstData = $wb.Worksheets.Item("2022-Q3")
$dstData.Range("A2:H9").PasteSpecial(-4122)      # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Fill in the "2022-Q3" fund-holding data.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-Q3")

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$rows = @(
    @{A=0; B="161903"; C="万家行业优选混合（LOF）";     D="83.48"; E="93.78"; F="3.95"; G="3.2975"; H=10},
    @{A=1; B="008120"; C="万家自主创新混合A";          D="22.86"; E="93.71"; F="5.26"; G="1.2024"; H=10},
    @{A=2; B="008633"; C="万家科技创新混合A";          D="2.44";  E="93.37"; F="7.95"; G="0.1940"; H=4},
    @{A=3; B="008634"; C="万家科技创新混合C";          D="1.71";  E="93.37"; F="7.95"; G="0.1359"; H=4},
    @{A=4; B="008121"; C="万家自主创新混合C";          D="2.24";  E="93.71"; F="5.26"; G="0.1178"; H=10},
    @{A=5; B="506008"; C="长城科创两年定开混合A";       D="3.12";  E="78.09"; F="3.00"; G="0.0936"; H=5},
    @{A=6; B="004223"; C="金信多策略精选灵活配置混合";   D="0.32";  E="92.79"; F="5.20"; G="0.0166"; H=6},
    @{A=7; B="012793"; C="长城科创两年定开混合C";       D="0.11";  E="78.09"; F="3.00"; G="0.0033"; H=5}
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A" + $r).Value = $row.A

    $ws.Range("B" + $r).Value = "'" + $row.B
    $ws.Range("B" + $r).Style = "Normal"

    $ws.Range("C" + $r).Value = $row.C

    $ws.Range("D" + $r).Value = "'" + $row.D
    $ws.Range("D" + $r).Style = "Normal"

    $ws.Range("E" + $r).Value = "'" + $row.E
    $ws.Range("E" + $r).Style = "Normal"

    $ws.Range("F" + $r).Value = "'" + $row.F
    $ws.Range("F" + $r).Style = "Normal"

    $ws.Range("G" + $r).Value = "'" + $row.G
    $ws.Range("G" + $r).Style = "Normal"

    $ws.Range("H" + $r).Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert a new row for 2022-Q3 above the
#    existing data (pushing the rest down by one row).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the data-row formatting (copied from the row directly below,
# which still carries the original data-row style) to the freshly inserted
# blank row.
$summary2 = $wb.Worksheets.Item("总计")
$summary2.Range("A3:D3").Copy()
$summary2.Range("A2:D2").PasteSpecial(-4122)     # xlPasteFormats

$summary3 = $wb.Worksheets.Item("总计")
$summary3.Range("A2").Value = 0
$summary3.Range("B2").Value = "2022-Q3"
$summary3.Range("C2").Value = 8
$summary3.Range("D2").Value = 5.06

# Renumber the 0-based index column (A) for the rows that shifted down.
for ($row = 3; $row -le 8; $row++) {
    $summary3.Cells.Item($row, 1).Value = $row - 2
}

Write-Output "2022-Q3 sheet added and 总计 summary updated"
